$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 111
$ws.Range("F2").Value = 80
$ws.Range("H2").Value = 86

$ws.Range("E3").Value = 46
$ws.Range("F3").Value = 34
$ws.Range("H3").Value = 35

$ws.Range("E5").Value = 154

$ws.Range("E7").Value = 42

$ws.Range("E9").Value = 12

$ws.Range("E10").Value = 670
$ws.Range("F10").Value = 374
$ws.Range("H10").Value = 470

$ws.Range("E11").Value = 439
$ws.Range("F11").Value = 248
$ws.Range("H11").Value = 312

$ws.Range("E12").Value = 672
$ws.Range("F12").Value = 406
$ws.Range("H12").Value = 492

$ws.Range("E13").Value = 157

$ws.Range("E14").Value = 137
$ws.Range("F14").Value = 81
$ws.Range("H14").Value = 115

$ws.Range("E15").Value = 195
$ws.Range("F15").Value = 92
$ws.Range("G15").Value = 51
$ws.Range("H15").Value = 143

$ws.Range("E16").Value = 229
$ws.Range("F16").Value = 133
$ws.Range("H16").Value = 181

$ws.Range("E17").Value = 124
$ws.Range("F17").Value = 68
$ws.Range("H17").Value = 92

$ws.Range("E18").Value = 59
$ws.Range("F18").Value = 31
$ws.Range("H18").Value = 48

$ws.Range("F20").Value = 43
$ws.Range("H20").Value = 80

$ws.Range("E21").Value = 149
$ws.Range("F21").Value = 89
$ws.Range("H21").Value = 120

$ws.Range("E22").Value = 192
$ws.Range("F22").Value = 107
$ws.Range("H22").Value = 149

$ws.Range("E23").Value = 222
$ws.Range("F23").Value = 114
$ws.Range("G23").Value = 52
$ws.Range("H23").Value = 166

$ws.Range("E24").Value = 252
$ws.Range("F24").Value = 151
$ws.Range("H24").Value = 181

$ws.Range("E25").Value = 315
$ws.Range("F25").Value = 178
$ws.Range("H25").Value = 238

$ws.Range("E26").Value = 183

$ws.Range("E27").Value = 372
$ws.Range("F27").Value = 205
$ws.Range("G27").Value = 82
$ws.Range("H27").Value = 287

$ws.Range("E28").Value = 221
$ws.Range("F28").Value = 110
$ws.Range("H28").Value = 162

$ws.Range("E29").Value = 191
$ws.Range("F29").Value = 116
$ws.Range("H29").Value = 157

$ws.Range("E30").Value = 248
$ws.Range("F30").Value = 155
$ws.Range("H30").Value = 207

$ws.Range("E32").Value = 208
$ws.Range("F32").Value = 133
$ws.Range("H32").Value = 171

$ws.Range("E33").Value = 321
$ws.Range("F33").Value = 179
$ws.Range("H33").Value = 269

$ws.Range("E34").Value = 246
$ws.Range("F34").Value = 173
$ws.Range("H34").Value = 211

$ws.Range("E35").Value = 175
$ws.Range("F35").Value = 124
$ws.Range("H35").Value = 151

$ws.Range("F36").Value = 58
$ws.Range("H36").Value = 68

$ws.Range("E37").Value = 189
$ws.Range("F37").Value = 107
$ws.Range("H37").Value = 144

$ws.Range("E38").Value = 103
$ws.Range("F38").Value = 63
$ws.Range("G38").Value = 17
$ws.Range("H38").Value = 80

$ws.Range("E40").Value = 296
$ws.Range("F40").Value = 153
$ws.Range("H40").Value = 233

$ws.Range("E41").Value = 427
$ws.Range("F41").Value = 214
$ws.Range("H41").Value = 306

$ws.Range("E42").Value = 437
$ws.Range("F42").Value = 255
$ws.Range("H42").Value = 316

$ws.Range("E43").Value = 139

$ws.Range("E44").Value = 352
$ws.Range("F44").Value = 188
$ws.Range("H44").Value = 256

$ws.Range("E45").Value = 177
$ws.Range("F45").Value = 99
$ws.Range("H45").Value = 138

$ws.Range("E46").Value = 382
$ws.Range("F46").Value = 223
$ws.Range("G46").Value = 64
$ws.Range("H46").Value = 287

$ws.Range("E47").Value = 526
$ws.Range("F47").Value = 298
$ws.Range("H47").Value = 390

$ws.Range("E48").Value = 257
$ws.Range("F48").Value = 125
$ws.Range("H48").Value = 169

$ws.Range("E49").Value = 331
$ws.Range("F49").Value = 167
$ws.Range("H49").Value = 254

$ws.Range("E50").Value = 270
$ws.Range("F50").Value = 146
$ws.Range("G50").Value = 73
$ws.Range("H50").Value = 219

$ws.Range("E51").Value = 262
$ws.Range("F51").Value = 130
$ws.Range("H51").Value = 204

$ws.Range("E52").Value = 32
$ws.Range("F52").Value = 15
$ws.Range("H52").Value = 23

